$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.251.59'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '2.421.57'
$ws.Range("E3").Value = '  +1.73%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.66%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +1.34%  '
$ws.Range("D9").Value = '2.419.91'
$ws.Range("E9").Value = '  +1.58%  '
$ws.Range("E10").Value = '  +0.62%  '
$ws.Range("E11").Value = '  -1.87%  '
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.89'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("D16").Value = '2.859.24'
$ws.Range("E16").Value = '  +1.72%  '
$ws.Range("D17").Value = '62.002.77'
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("D18").Value = '2.420.93'
$ws.Range("E18").Value = '  +1.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.79'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.53%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.60'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.99%  '
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("E26").Value = '  +2.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '582.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +12.11%  '
$ws.Range("D28").Value = '2.541.13'
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").Value = '0.0₃0945'
$ws.Range("E30").Value = '  +4.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.51%  '
$ws.Range("E33").Value = '  +1.17%  '
$ws.Range("E34").Value = '  +1.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.46%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.22%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.80'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.385'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '152.54'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.66'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.72%  '
$ws.Range("E42").Value = '  -2.78%  '
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("E44").Value = '  +7.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '150.39'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.53%  '
$ws.Range("E46").Value = '  +1.41%  '
$ws.Range("E47").Value = '  +2.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.26'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.595'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.05%  '
$ws.Range("E50").Value = '  +1.77%  '
$ws.Range("E51").Value = '  +1.70%  '
